$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Q" (2020) column to the table. Copying the existing "P" (2019)
# column first brings over the matching number format / font / borders for
# each row, then the actual 2020 figures are written on top of that.
$ws.Range("P4:P14").Copy($ws.Range("Q4:Q14"))

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 4.5999999999999996
$ws.Range("Q6").Value = 4.2
$ws.Range("Q7").Value = 1.3
$ws.Range("Q8").Value = 10.8
$ws.Range("Q9").Value = 6.5
$ws.Range("Q10").Value = 2.9
$ws.Range("Q11").Value = 2.6
$ws.Range("Q12").Value = 13.1
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 1.3

# Reproduce the author's final UI selection state: whole column T selected.
$ws.Range("T:T").Select()
